$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet from "Session" to "Anatomy"
$ws.Name = "Anatomy"

# Append a new QR-scan log row (row 24) mirroring the existing rows' layout.
# Column A holds a numeric-looking Student ID that must stay text (like the
# rows above it), so mark the cell as Text before writing the value and then
# clear the number-format styling we only used to force the text type.
$ws.Cells.Item(24, 1).NumberFormat = "@"
$ws.Cells.Item(24, 1).Value = "234035"
$ws.Cells.Item(24, 1).ClearFormats()

$ws.Cells.Item(24, 2).Value = "Anatomy"
$ws.Cells.Item(24, 3).Value = "15/10/2025"
$ws.Cells.Item(24, 4).Value = "12:54:40"
$ws.Cells.Item(24, 5).Value = "Scan"
$ws.Cells.Item(24, 6).Value = "hananragab@med.asu.edu.eg"
